$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 15.20717213437361
$ws.Cells.Item(2, 3).Value = 11.23098386416958
$ws.Cells.Item(2, 4).Value = 3.691276806880144
$ws.Cells.Item(2, 5).Value = 16.58230678704313
$ws.Cells.Item(2, 6).Value = 19.73925079067827
$ws.Cells.Item(2, 15).Value = 17.46919510004873

$ws.Cells.Item(3, 2).Value = 14.3869011330943
$ws.Cells.Item(3, 3).Value = 10.61409342803636
$ws.Cells.Item(3, 4).Value = 3.65097749453512
$ws.Cells.Item(3, 5).Value = 15.63480355221651
$ws.Cells.Item(3, 6).Value = 19.73891381041926
$ws.Cells.Item(3, 15).Value = 17.54996281088463

$ws.Cells.Item(4, 2).Value = 13.85872582943154
$ws.Cells.Item(4, 3).Value = 10.21494742836907
$ws.Cells.Item(4, 4).Value = 3.626050790119419
$ws.Cells.Item(4, 5).Value = 15.02787072369468
$ws.Cells.Item(4, 6).Value = 19.74890707209335
$ws.Cells.Item(4, 15).Value = 17.6074872766523

$ws.Cells.Item(5, 2).Value = 13.63750985203694
$ws.Cells.Item(5, 3).Value = 10.04725985989433
$ws.Cells.Item(5, 4).Value = 3.615854635603046
$ws.Cells.Item(5, 5).Value = 14.77448087684881
$ws.Cells.Item(5, 6).Value = 19.75552969658342
$ws.Cells.Item(5, 15).Value = 17.63290483321494

$ws.Cells.Item(6, 2).Value = 13.60042211023479
$ws.Cells.Item(6, 3).Value = 10.01911457157782
$ws.Cells.Item(6, 4).Value = 3.61415948184738
$ws.Cells.Item(6, 5).Value = 14.73204855493091
$ws.Cells.Item(6, 6).Value = 19.75678291737364
$ws.Cells.Item(6, 15).Value = 17.63724417227236

$ws.Cells.Item(7, 2).Value = 13.8557663731598
$ws.Cells.Item(7, 3).Value = 10.21270617651269
$ws.Cells.Item(7, 4).Value = 3.625913426264311
$ws.Cells.Item(7, 5).Value = 15.02447756294644
$ws.Cells.Item(7, 6).Value = 19.7489860829959
$ws.Cells.Item(7, 15).Value = 17.60782209203203

$ws.Cells.Item(8, 2).Value = 14.92953414443551
$ws.Cells.Item(8, 3).Value = 11.02255680238469
$ws.Cells.Item(8, 4).Value = 3.677424094981955
$ws.Cells.Item(8, 5).Value = 16.26097475304845
$ws.Cells.Item(8, 6).Value = 19.73700978143365
$ws.Cells.Item(8, 15).Value = 17.49538620089033

$ws.Cells.Item(9, 2).Value = 16.83407210579252
$ws.Cells.Item(9, 3).Value = 12.44595959867729
$ws.Cells.Item(9, 4).Value = 3.776655678919636
$ws.Cells.Item(9, 5).Value = 18.59159515867276
$ws.Cells.Item(9, 6).Value = 19.7949446420826
$ws.Cells.Item(9, 15).Value = 17.33865376030681

$ws.Cells.Item(10, 2).Value = 18.1039505849988
$ws.Cells.Item(10, 3).Value = 13.38825905305182
$ws.Cells.Item(10, 4).Value = 3.848041828199405
$ws.Cells.Item(10, 5).Value = 20.24904979014259
$ws.Cells.Item(10, 6).Value = 19.8875623629176
$ws.Cells.Item(10, 15).Value = 17.26342360399078

$ws.Cells.Item(11, 2).Value = 18.65257390207194
$ws.Cells.Item(11, 3).Value = 13.79405483640517
$ws.Cells.Item(11, 4).Value = 3.880098268908473
$ws.Cells.Item(11, 5).Value = 20.96075776667514
$ws.Cells.Item(11, 6).Value = 19.94057624441709
$ws.Cells.Item(11, 15).Value = 17.23807508470405

$ws.Cells.Item(12, 2).Value = 18.85608603829845
$ws.Cells.Item(12, 3).Value = 13.9444083126921
$ws.Cells.Item(12, 4).Value = 3.892169936680048
$ws.Cells.Item(12, 5).Value = 21.22422159318808
$ws.Cells.Item(12, 6).Value = 19.96221098400179
$ws.Cells.Item(12, 15).Value = 17.22976779147638

$ws.Cells.Item(13, 2).Value = 18.81244545916545
$ws.Cells.Item(13, 3).Value = 13.91217465645377
$ws.Cells.Item(13, 4).Value = 3.889573200256137
$ws.Cells.Item(13, 5).Value = 21.16774822534621
$ws.Cells.Item(13, 6).Value = 19.95748231843642
$ws.Cells.Item(13, 15).Value = 17.23149923406023

$ws.Cells.Item(14, 2).Value = 18.66940224276311
$ws.Cells.Item(14, 3).Value = 13.80649101585829
$ws.Cells.Item(14, 4).Value = 3.881092804767366
$ws.Cells.Item(14, 5).Value = 20.98255403741254
$ws.Cells.Item(14, 6).Value = 19.94232493546568
$ws.Cells.Item(14, 15).Value = 17.23736566723153

$ws.Cells.Item(15, 2).Value = 18.58123041413059
$ws.Cells.Item(15, 3).Value = 13.74132471131175
$ws.Cells.Item(15, 4).Value = 3.875889323592601
$ws.Cells.Item(15, 5).Value = 20.86833110078487
$ws.Cells.Item(15, 6).Value = 19.93324347833317
$ws.Cells.Item(15, 15).Value = 17.24112767902822

$ws.Cells.Item(16, 2).Value = 18.06750654776319
$ws.Cells.Item(16, 3).Value = 13.36127700406613
$ws.Cells.Item(16, 4).Value = 3.845937828067923
$ws.Cells.Item(16, 5).Value = 20.20169026644631
$ws.Cells.Item(16, 6).Value = 19.88431638632261
$ws.Cells.Item(16, 15).Value = 17.2652601142246

$ws.Cells.Item(17, 2).Value = 17.74486417280591
$ws.Cells.Item(17, 3).Value = 13.12225453602907
$ws.Cells.Item(17, 4).Value = 3.827451312355486
$ws.Cells.Item(17, 5).Value = 19.78192449817088
$ws.Cells.Item(17, 6).Value = 19.85708572443544
$ws.Cells.Item(17, 15).Value = 17.28234886940117

$ws.Cells.Item(18, 2).Value = 17.55655712956199
$ws.Cells.Item(18, 3).Value = 12.98262434868256
$ws.Cells.Item(18, 4).Value = 3.816779490824612
$ws.Cells.Item(18, 5).Value = 19.53650500541305
$ws.Cells.Item(18, 6).Value = 19.84244815258743
$ws.Cells.Item(18, 15).Value = 17.29301227593966

$ws.Cells.Item(19, 2).Value = 17.49233256013565
$ws.Cells.Item(19, 3).Value = 12.93497945426967
$ws.Cells.Item(19, 4).Value = 3.813159748368181
$ws.Cells.Item(19, 5).Value = 19.45272532330779
$ws.Cells.Item(19, 6).Value = 19.83766822721991
$ws.Cells.Item(19, 15).Value = 17.29676553953992

$ws.Cells.Item(20, 2).Value = 17.77949321805978
$ws.Cells.Item(20, 3).Value = 13.14792168900441
$ws.Cells.Item(20, 4).Value = 3.829423310232393
$ws.Cells.Item(20, 5).Value = 19.8270210794926
$ws.Cells.Item(20, 6).Value = 19.85987843969546
$ws.Cells.Item(20, 15).Value = 17.28044326825091

$ws.Cells.Item(21, 2).Value = 18.71153299112191
$ws.Cells.Item(21, 3).Value = 13.8376229314562
$ws.Cells.Item(21, 4).Value = 3.883585589698038
$ws.Cells.Item(21, 5).Value = 21.03711383061446
$ws.Cells.Item(21, 6).Value = 19.94673476135729
$ws.Cells.Item(21, 15).Value = 17.23560738050466

$ws.Cells.Item(22, 2).Value = 19.29594365801019
$ws.Cells.Item(22, 3).Value = 14.26906731943327
$ws.Cells.Item(22, 4).Value = 3.918587586866277
$ws.Cells.Item(22, 5).Value = 21.79276318518961
$ws.Cells.Item(22, 6).Value = 20.01258546668178
$ws.Cells.Item(22, 15).Value = 17.21383959900516

$ws.Cells.Item(23, 2).Value = 18.98631264627127
$ws.Cells.Item(23, 3).Value = 14.04057124366009
$ws.Cells.Item(23, 4).Value = 3.899944992668285
$ws.Cells.Item(23, 5).Value = 21.39266975740218
$ws.Cells.Item(23, 6).Value = 19.97661112012243
$ws.Cells.Item(23, 15).Value = 17.22476317425515

$ws.Cells.Item(24, 2).Value = 17.76384619196211
$ws.Cells.Item(24, 3).Value = 13.13632446018795
$ws.Cells.Item(24, 4).Value = 3.828531905456983
$ws.Cells.Item(24, 5).Value = 19.80664565929584
$ws.Cells.Item(24, 6).Value = 19.85861268274044
$ws.Cells.Item(24, 15).Value = 17.28130217875398

$ws.Cells.Item(25, 2).Value = 16.34116850601542
$ws.Cells.Item(25, 3).Value = 12.07888989524571
$ws.Cells.Item(25, 4).Value = 3.75004685495666
$ws.Cells.Item(25, 5).Value = 17.94350157946699
$ws.Cells.Item(25, 6).Value = 19.77049761042823
$ws.Cells.Item(25, 15).Value = 17.37411645601009
